# Update "Pais" sheet with refreshed COVID case counts and the new
# "datos actualizados" timestamp.
#
# Two country pairs swap rank order because of the updated totals:
#   - Austria overtakes Argelia (rows 64/65)
#   - Eslovaquia overtakes Madagascar (rows 92/93)
# so column A (country name) is rewritten for those four rows as well as
# the numeric columns B:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header (A1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 11:51"

function Set-Row($Row, $Name, $B, $C, $D, $E, $F, $G, $H) {
    if ($Name -ne $null) {
        $ws.Cells.Item($Row, 1).Value = $Name
    }
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
}

# Rows whose totals changed but whose rank/order stayed the same
Set-Row 4   $null 7834289 526  5025910 2590629 0 12 217750
Set-Row 19  $null 375870  1278 289912  80481   0 17 5477
Set-Row 27  $null 286646  1310 225189  59578   0 15 1879
Set-Row 39  $null 116338  4739 77875   35544   0 52 2919
Set-Row 91  $null 19446   457  16695   2438    0 3  313
Set-Row 102 $null 11580   235  8500    2734    0 0  346
Set-Row 136 $null 4488    0    3296    1179    0 0  13
Set-Row 140 $null 3809    49   2906    835     0 1  68
Set-Row 188 $null 282     1    277     5       0 0  0
Set-Row 196 $null 137     2    117     19      0 0  1

# Austria now ranks above Argelia -> rows 64/65 swap contents
Set-Row 64 "Austria" 53188 1131 42039 10307 0 4 842
Set-Row 65 "Argelia" 52658 0    36958 13917 0 0 1783

# Eslovaquia now ranks above Madagascar -> rows 92/93 swap contents
Set-Row 92 "Eslovaquia" 16910 1184 5452  11401 0 0 57
Set-Row 93 "Madagascar" 16654 0    15910 509   0 0 235
